$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sectors")
$ws.Activate()

# Rows 33-36 (1B2_Fugitive-petr, 1B2b_Fugitive-NG-prod, 1B2b_Fugitive-NG-distr,
# 1B2d_Fugitive-other-energy) switch their activity driver from
# "refinery-and-natural-gas" to "pop", and their units from the text "kt" to
# the numeric value 1000 - matching the other "pop"-driven sector rows
# (oil production driver / extension updates).
$ws.Range("C33:C36").Value = "pop"
$ws.Range("D33:D36").Value = 1000

# Column C widened slightly to comfortably fit the "pop" activity values.
$ws.Columns.Item(3).ColumnWidth = 14.83

# Leave the view scrolled/selected over the edited range, as happened while
# reviewing the change.
$ws.Range("C34:E36").Select()

$wb.Save()
